# Apply updated crypto price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.877.25'
$ws.Range('E2').Value = '  +5.20%  '
$ws.Range('D3').Value = '3.115.61'
$ws.Range('E3').Value = '  +3.46%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '''583.62'
$ws.Range('E5').Value = '  +3.19%  '
$ws.Range('D6').Value = '''145.05'
$ws.Range('E6').Value = '  +3.15%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '3.107.51'
$ws.Range('E8').Value = '  +3.57%  '
$ws.Range('D9').Value = '''0.529'
$ws.Range('E9').Value = '  +1.62%  '
$ws.Range('E10').Value = '  +11.20%  '
$ws.Range('E11').Value = '  +7.33%  '
$ws.Range('D12').Value = '''0.468'
$ws.Range('E12').Value = '  +1.69%  '
$ws.Range('E13').Value = '  +6.82%  '
$ws.Range('D14').Value = '''35.61'
$ws.Range('E14').Value = '  +4.84%  '
$ws.Range('D16').Value = '3.630.71'
$ws.Range('E16').Value = '  +3.21%  '
$ws.Range('E17').Value = '  -0.37%  '
$ws.Range('D18').Value = '3.109.49'
$ws.Range('E18').Value = '  +3.12%  '
$ws.Range('D19').Value = '62.785.28'
$ws.Range('E19').Value = '  +5.06%  '
$ws.Range('D20').Value = '''467.59'
$ws.Range('E20').Value = '  +6.80%  '
$ws.Range('D21').Value = '''14.08'
$ws.Range('E21').Value = '  +3.03%  '
$ws.Range('D22').Value = '''0.728'
$ws.Range('E22').Value = '  +1.43%  '
$ws.Range('D23').Value = '''7.55'
$ws.Range('E23').Value = '  +6.25%  '
$ws.Range('D24').Value = '''13.32'
$ws.Range('E24').Value = '  -0.66%  '
$ws.Range('D25').Value = '''82.02'
$ws.Range('E25').Value = '  +1.71%  '
$ws.Range('E26').Value = '  -0.13%  '
$ws.Range('E27').Value = '  +0.99%  '
$ws.Range('E28').Value = '  +4.83%  '
$ws.Range('B29').Value = 'FirstDigitalUSD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D29').Value = '''1.00'
$ws.Range('E29').Value = '  -0.09%  '
$ws.Range('B30').Value = 'RenderToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D30').Value = '''8.27'
$ws.Range('E30').Value = '  +5.37%  '
$ws.Range('D31').Value = '''6.82'
$ws.Range('E31').Value = '  +7.58%  '
$ws.Range('D32').Value = '''27.00'
$ws.Range('E32').Value = '  +3.95%  '
$ws.Range('D33').Value = '''0.111'
$ws.Range('E33').Value = '  +3.46%  '
$ws.Range('D34').Value = '0.0₃0870'
$ws.Range('E34').Value = '  +10.62%  '
$ws.Range('E35').Value = '  +12.75%  '
$ws.Range('E36').Value = '  +4.02%  '
$ws.Range('E37').Value = '  +1.92%  '
$ws.Range('D38').Value = '''3.27'
$ws.Range('E38').Value = '  +16.68%  '
$ws.Range('D39').Value = '''50.95'
$ws.Range('E39').Value = '  +3.37%  '
$ws.Range('D40').Value = '''434.42'
$ws.Range('E40').Value = '  +7.75%  '
$ws.Range('E41').Value = '  +1.89%  '
$ws.Range('D42').Value = '2.931.58'
$ws.Range('E42').Value = '  +5.94%  '
$ws.Range('E43').Value = '  +4.52%  '
$ws.Range('D44').Value = '''0.279'
$ws.Range('E44').Value = '  +10.11%  '
$ws.Range('E45').Value = '  +3.50%  '
$ws.Range('E46').Value = '  +5.50%  '
$ws.Range('D47').Value = '''35.72'
$ws.Range('E47').Value = '  +5.29%  '
$ws.Range('D48').Value = '''0.999'
$ws.Range('E48').Value = '  -0.01%  '
$ws.Range('D49').Value = '''123.61'
$ws.Range('E49').Value = '  +0.03%  '
$ws.Range('E50').Value = '  +0.63%  '
$ws.Range('E51').Value = '  +4.03%  '
